$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row before "Description" (currently row 11) to make room for "Jurisdiction"
$ws.Rows.Item(11).Insert()

# Fill the new Jurisdiction row
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update Version value (row 3)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"
